# Update the "取得日時" (retrieved datetime) column for all existing data
# rows on the "ランサーズ" sheet to reflect the new scrape timestamp.
#
# Commit message: Append: 2025-10-31 01:18 JST
# Diff: A2:A13 change from "2025-10-30 18:35:10" to "2025-10-31 01:18:52"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-31 01:18:52"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
